$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("A9").Value = "For India, the water output unit is cubic km, which is equivalent to Tl (teraliters, or 10^12 liters)"
